# Fixed clients class, so it works with % instead of decimals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set number format to Text (numFmtId 49) for the discount columns (C2:E6),
# matching the change in cellXfs from numFmtId 2 to numFmtId 49.
$ws.Range("C2:E6").NumberFormat = "@"

# Column-major order (C2:C6, then D2:D6, then E2:E6) so that new shared
# strings get registered in the same order as the target workbook.
$colC = @("5%", "4%", "3%", "2%", "0%")
$colD = @("0%", "1%", "1%", "3%", "5%")
$colE = @("2%", "2%", "3%", "5%", "7%")

for ($i = 0; $i -lt 5; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $colC[$i]
}
for ($i = 0; $i -lt 5; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $colD[$i]
}
for ($i = 0; $i -lt 5; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = $colE[$i]
}

# Update the selected cell/range as in the diff
$ws.Range("C15").Select()

$wb.Save()
